$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns K, L, M, N, P on row 2 become numeric 1
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1
$ws.Range("N2").Value = 1
$ws.Range("P2").Value = 1

# O2 becomes text date string (keep as literal text, not an Excel date serial)
$ws.Range("O2").NumberFormat = "@"
$ws.Range("O2").Value = "2025-02-28"
$ws.Range("O2").Style = "Normal"

# Q2:V2 become the check mark emoji
$ws.Range("Q2").Value = "✅"
$ws.Range("R2").Value = "✅"
$ws.Range("S2").Value = "✅"
$ws.Range("T2").Value = "✅"
$ws.Range("U2").Value = "✅"
$ws.Range("V2").Value = "✅"
